$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Bruno Raniere"
$ws.Range("B3").Value = "+55 11 91234-5678"
$ws.Range("C3").Value = "Oi, como posso usar o WhatsApp GPT?"
$ws.Range("D3").Value = "13ct49b3764trc76134bc613rtb79346tcr973bcr2363"
